$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test-case step name in B1: readProperties -> setProperty
$ws.Range("B1").Value = "setProperty"

# Update the JSON payload in B3 for the new step
$ws.Range("B3").Value = '{"name":"hugang","admin":"administrator"}'

# The longer wrapped text now needs two lines, so the row grows taller
$ws.Rows.Item(3).RowHeight = 28.8

# Cursor/selection ends on E9, matching the saved file's sheetView selection
$ws.Range("E9").Select()
